# Add new register rows (10-16) to the worksheet, matching the data
# appended in the source workbook update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows data: id (number), nome, address, postal_code (all stored as text)
$newRows = @(
    @{ Row = 10; Id = 691; Nome = "Vinicius Ferreir";  Address = "Av Leoncio de Magalhaes, 699"; Postal = "2042010" },
    @{ Row = 11; Id = 578; Nome = "Vinicius Ferreira"; Address = "Av Leoncio de Magalhaes, 699"; Postal = "2042010" },
    @{ Row = 12; Id = 130; Nome = "Duda Santos";       Address = "Rua Alvares Machado";          Postal = "264795" },
    @{ Row = 13; Id = 186; Nome = "Maria Oliveira";    Address = "Rua Alvares 2";                Postal = "2313453" },
    @{ Row = 14; Id = 21;  Nome = "Maria Oliveira";    Address = "Rua Alvares 2";                Postal = "2313453" },
    @{ Row = 15; Id = 654; Nome = "Maria Santos";      Address = "Rua jorge amado";              Postal = "91701153890" },
    @{ Row = 16; Id = 757; Nome = "Maria Santos";      Address = "Rua jorge amado";              Postal = "91701153890" }
)

# Format the postal_code column for the new rows as Text up front so the
# numeric-looking codes are written out as strings (matching column D of
# rows 7-9 which are also stored as text), using a single shared style.
$ws.Range("D10:D16").NumberFormat = "@"

foreach ($r in $newRows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.Id
    $ws.Cells.Item($r.Row, 2).Value = $r.Nome
    $ws.Cells.Item($r.Row, 3).Value = $r.Address
    $ws.Cells.Item($r.Row, 4).Value = $r.Postal
}
